$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values (rows 2 and 3)
# Row 2
$ws.Range("D2").Value = 0.021
$ws.Range("E2").Value = -0.141
$ws.Range("G2").Value = 0.09460431654676257
$ws.Range("H2").Value = 0.09460431654676257
$ws.Range("I2").Value = 0.06750599520383692
$ws.Range("J2").Value = 0.06416067146282972
$ws.Range("K2").Value = 5.31
$ws.Range("L2").Value = 0.06366906474820143
$ws.Range("M2").Value = 5.366699999999999
$ws.Range("N2").Value = 0.04032081141998497
$ws.Range("O2").Value = 1.010677966101695
$ws.Range("P2").Value = 5.366699999999999
$ws.Range("Q2").Value = 0.04032081141998497
$ws.Range("R2").Value = 1.010677966101695
$ws.Range("U2").Value = 94
$ws.Range("V2").Value = 0.7062359128474831
$ws.Range("W2").Value = 0.02673716012084592
$ws.Range("X2").Value = 0.06839863771989911
$ws.Range("Y2").Value = -0.04166147759905319
$ws.Range("Z2").Value = 0.6430223592906709
$ws.Range("AA2").Value = 0.04125674633770239
$ws.Range("AB2").Value = 0.06839863771989911
$ws.Range("AC2").Value = -0.02714189138219672
$ws.Range("AG2").Value = -94
$ws.Range("AJ2").Value = -2.404092071611254
$ws.Range("AK2").Value = -0.9523809523809524
$ws.Range("AP2").Value = -15.98639455782313

# Row 3
$ws.Range("D3").Value = 0.021
$ws.Range("E3").Value = -0.141
$ws.Range("G3").Value = 0.09460431654676257
$ws.Range("H3").Value = 0.09460431654676257
$ws.Range("I3").Value = 0.06750599520383692
$ws.Range("J3").Value = 0.06416067146282972
$ws.Range("K3").Value = 5.31
$ws.Range("L3").Value = 0.06366906474820143
$ws.Range("M3").Value = 5.366699999999999
$ws.Range("N3").Value = 0.04032081141998497
$ws.Range("O3").Value = 1.010677966101695
$ws.Range("P3").Value = 5.366699999999999
$ws.Range("Q3").Value = 0.04032081141998497
$ws.Range("R3").Value = 1.010677966101695
$ws.Range("U3").Value = 94
$ws.Range("V3").Value = 0.7062359128474831
$ws.Range("W3").Value = 0.02673716012084592
$ws.Range("X3").Value = 0.06839863771989911
$ws.Range("Y3").Value = -0.04166147759905319
$ws.Range("Z3").Value = 0.6430223592906709
$ws.Range("AA3").Value = 0.04125674633770239
$ws.Range("AB3").Value = 0.06839863771989911
$ws.Range("AC3").Value = -0.02714189138219672
$ws.Range("AG3").Value = -94
$ws.Range("AJ3").Value = -2.404092071611254
$ws.Range("AK3").Value = -0.9523809523809524
$ws.Range("AP3").Value = -15.98639455782313

